$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per data row (Excel row 2-25), columns B,C,D,F,G,J,K,L
# Row index (array element 0) corresponds to Excel row 2, etc.
$data = @(
    @(1.806795917023464, 0.01365998443658611, 0.0216406295584477, 6.990267896841829, 0.00269604581964237, 0.2890488050483313, 1.232829515477675, 0.285116032269876),
    @(1.80015809124103, 0.01177974038080976, 0.01910724806568709, 6.796685944482363, 0.002701233681433745, 0.285473506385415, 1.223749541883478, 0.2874226506143316),
    @(1.797453289979302, 0.01063189762125205, 0.01755943812875671, 6.678346016927122, 0.00270458450299471, 0.2833088728896698, 1.219257744936954, 0.2890464724699129),
    @(1.7966956522163, 0.01016570710898179, 0.01693048644654738, 6.630248799866166, 0.002705991739570582, 0.2824344081125716, 1.217699403914722, 0.2897604115504535),
    @(1.796590656294569, 0.01008838802034973, 0.01682615367209195, 6.622269877384838, 0.002706227935967844, 0.2822896637425032, 1.217457068064761, 0.2898821153364608),
    @(1.797441677138409, 0.01062560421432579, 0.01755094879052166, 6.677696850445358, 0.002704603312369826, 0.2832970486865207, 1.219235627285485, 0.2890558894337545),
    @(1.804222587960396, 0.01301023618433561, 0.02076543609400971, 6.923410682356177, 0.00269780034202141, 0.2878096522070663, 1.229473692007303, 0.2858682937321504),
    @(1.828407489318636, 0.01774403201945063, 0.02713675687495254, 7.409566419040345, 0.002685765909649445, 0.2969048711120621, 1.258164833851083, 0.2812632930246792),
    @(1.852835275939498, 0.02126413045307629, 0.03186862192069384, 7.769664403187704, 0.002677711191366578, 0.303742162007886, 1.284527678821746, 0.2788825699839421),
    @(1.865399553513413, 0.02287623683641016, 0.03403440093862287, 7.934181529898467, 0.002674215789434465, 0.3068874339571579, 1.297675631128243, 0.2780171197808272),
    @(1.870366465655763, 0.02348837417783045, 0.03485658938710401, 7.996586098557373, 0.002672916283265264, 0.3080835715371961, 1.302821094008237, 0.2777206711609068),
    @(1.869287447742863, 0.02335646343173892, 0.03467942256381207, 7.983141403632771, 0.00267319508440807, 0.3078257344225506, 1.301705508276058, 0.2777831256102985),
    @(1.865803992022364, 0.02292656365766277, 0.03410200081629, 7.939313454515911, 0.002674108395680686, 0.3069857385007708, 1.298095609360729, 0.2779921039272892),
    @(1.863697514228221, 0.02266345825037774, 0.03374858547803683, 7.912481425613805, 0.002674670963181084, 0.3064718817370817, 1.295906155241738, 0.2781241824592229),
    @(1.852043413390248, 0.02115900256809766, 0.03172736240473739, 7.758927373656036, 0.002677943007642458, 0.3035373199864964, 1.283691718587903, 0.2789435064082895),
    @(1.845266090149352, 0.02023891438608416, 0.03049090424968881, 7.664910692444153, 0.002679993419304911, 0.3017460598920536, 1.276494808316301, 0.2795018511384981),
    @(1.84150459859697, 0.01971071319655948, 0.02978096231292682, 7.610901071931067, 0.002681188651456362, 0.3007190605502998, 1.272464054273087, 0.2798434735224475),
    @(1.840254482248952, 0.01953204319232782, 0.02954079583872726, 7.592625610191277, 0.002681596070014719, 0.3003718986454516, 1.271117968388864, 0.2799626581751724),
    @(1.845973403384789, 0.02033675406440238, 0.03062239829580449, 7.67491204024418, 0.00267977350598342, 0.3019364020180504, 1.277249676374026, 0.2794402951194357),
    @(1.866821490167212, 0.02305278945785005, 0.03427154659006248, 7.952183891131426, 0.002673839480155151, 0.3072323269638275, 1.299151398285062, 0.277929873068345),
    @(1.881665828738676, 0.02483765030657992, 0.03666852509454088, 8.134013860771859, 0.002670101816051039, 0.3107232273031002, 1.314436856294634, 0.2771250396882579),
    @(1.873631490106078, 0.02388410368268978, 0.03538805944829448, 8.036910055218357, 0.002672083861119776, 0.3088573281226701, 1.306189671923164, 0.2775379137309031),
    @(1.845653206949436, 0.02029251835757861, 0.03056294700256501, 7.670390300909901, 0.002679872877567051, 0.3018503395665633, 1.276908067747883, 0.2794680603582549),
    @(1.820697062064795, 0.01645652523781393, 0.02540490902737247, 7.277556917269266, 0.00268888267128258, 0.2944175166581218, 1.249477658322692, 0.2823329564005377)
)

$cols = @("B","C","D","F","G","J","K","L")
$startRow = 2

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $cellRef = "{0}{1}" -f $cols[$j], $rowNum
        $ws.Range($cellRef).Value = $rowValues[$j]
    }
}

Write-Host "Updated $($data.Count) rows"
